$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Name -> Geology
$ws.Range("B4").Value = "Geology"
$ws.Range("C4").Value = "Geology"

# Row 8: Ativação date (force text so Excel doesn't convert it to a date serial)
$ws.Range("B8").Value = "'01/01/2022"
$ws.Range("C8").Value = "'01/01/2022"

# Row 11: Objectives (English)
$ws.Range("B11").Value = "Provide basic knowledge about terrestrial materials and the main geological processes."
$ws.Range("C11").Value = "Provide basic knowledge about terrestrial materials and the main geological processes."

# Row 14: Programa resumido (Portuguese short syllabus)
$ws.Range("B14").Value = "Processos endógenos e exógenos da Terra. Materiais constituintes da crosta terrestre (minerais e rochas)."
$ws.Range("C14").Value = "Processos endógenos e exógenos da Terra. Materiais constituintes da crosta terrestre (minerais e rochas)."

# Row 15: Short syllabus (English)
$ws.Range("B15").Value = "Endogenous and exogenous processes of the Earth. Materials constituting the earth's crust (minerals and rocks)."
$ws.Range("C15").Value = "Endogenous and exogenous processes of the Earth. Materials constituting the earth's crust (minerals and rocks)."

# Row 16: Programa (Portuguese full syllabus)
$ws.Range("B16").Value = "Breve história da Geologia. Materiais constituintes da crosta terrestre (minerais e rochas). Origem e constituição do universo, do sistema solar e da Terra. Estrutura interna da Terra. Composição da Terra. Processos endógenos e exógenos (dinâmica interna e externa da Terra).  Teoria da tectônica de placas.  Rochas ígneas e vulcanismo. Rochas metamórficas e metamorfismo. Rochas sedimentares. Intemperismo, erosão, transporte de sedimentos.  Estrutura geológicas. Tempo geológico e estratigrafia."
$ws.Range("C16").Value = "Breve história da Geologia. Materiais constituintes da crosta terrestre (minerais e rochas). Origem e constituição do universo, do sistema solar e da Terra. Estrutura interna da Terra. Composição da Terra. Processos endógenos e exógenos (dinâmica interna e externa da Terra).  Teoria da tectônica de placas.  Rochas ígneas e vulcanismo. Rochas metamórficas e metamorfismo. Rochas sedimentares. Intemperismo, erosão, transporte de sedimentos.  Estrutura geológicas. Tempo geológico e estratigrafia."

# Row 17: Syllabus (English)
$ws.Range("B17").Value = "Brief history of geology. Materials constituting the earth's crust (minerals and rocks). Origin and constitution of the universe, the solar system and the earth. Internal structure of the earth. Composition of the earth. Endogenous and exogenous processes (internal and external dynamics of the earth). Plate tectonics theory.  Igneous rocks and vulcanismo. Metamorphic rocks and metamorphism. Sedimentary rocks. Weathering, erosion, sediment transport. Geological structure. Geological time and stratigraphy."
$ws.Range("C17").Value = "Brief history of geology. Materials constituting the earth's crust (minerals and rocks). Origin and constitution of the universe, the solar system and the earth. Internal structure of the earth. Composition of the earth. Endogenous and exogenous processes (internal and external dynamics of the earth). Plate tectonics theory.  Igneous rocks and vulcanismo. Metamorphic rocks and metamorphism. Sedimentary rocks. Weathering, erosion, sediment transport. Geological structure. Geological time and stratigraphy."

# Row 19: Método
$ws.Range("B19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
$ws.Range("C19").Value = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."

# Row 20: Critério
$ws.Range("B20").Value = "Média ponderada de provas  e atividades."
$ws.Range("C20").Value = "Média ponderada de provas  e atividades."

# Row 21: Norma de recuperação
$ws.Range("B21").Value = "1 (uma) prova escrita"
$ws.Range("C21").Value = "1 (uma) prova escrita"

# Row 22: Bibliografia
$ws.Range("B22").Value = "Bibliografia básica:PRESS, F.; SIEVER, R.; GROTZINGER, J.; JORDAN, T. H. Para entender a Terra. Porto Alegre: Bookman, 2008. 656p.REED, W.; MONROE, J. S. Fundamentos de Geologia. São Paulo: Cengage Learning, 2011. 508p.Bibliografia complementar:TEIXEIRA, W.; FAIRCHILD, T. R.; DE TOLEDO, M. C. M.; TAIOLI, F. Decifrando a Terra. São Paulo: Companhia Editora Nacional, 2003. 623p."
$ws.Range("C22").Value = "Bibliografia básica:PRESS, F.; SIEVER, R.; GROTZINGER, J.; JORDAN, T. H. Para entender a Terra. Porto Alegre: Bookman, 2008. 656p.REED, W.; MONROE, J. S. Fundamentos de Geologia. São Paulo: Cengage Learning, 2011. 508p.Bibliografia complementar:TEIXEIRA, W.; FAIRCHILD, T. R.; DE TOLEDO, M. C. M.; TAIOLI, F. Decifrando a Terra. São Paulo: Companhia Editora Nacional, 2003. 623p."
